{"js": "// Insert \", Data Fixing\" right after \"Bug Fixing\" in the bullet describing\n// systems enhancements / bug fixing work, turning:\n//   \"... Worked also on Bug Fixing and provided UAT and Production support.\"\n// into:\n//   \"... Worked also on Bug Fixing, Data Fixing and provided UAT and Production support.\"\n\nconst results = context.document.body.search(\"Bug Fixing\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \"Bug Fixing\" not found');\n}\n\n// Insert the new text immediately after the matched \"Bug Fixing\" phrase.\nconst hit = results.items[0];\nhit.insertText(\", Data Fixing\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert \", Data Fixing\" right after \"Bug Fixing\" in the bullet describing\n# systems enhancements / bug fixing work, turning:\n#   \"... Worked also on Bug Fixing and provided UAT and Production support.\"\n# into:\n#   \"... Worked also on Bug Fixing, Data Fixing and provided UAT and Production support.\"\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"Bug Fixing\", $true)\n\nif (-not $found) {\n    throw 'Target text \"Bug Fixing\" not found'\n}\n\n# $range now spans just the matched text (\"Bug Fixing\"); collapse to its end\n# so the new text is inserted immediately after it.\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertAfter(\", Data Fixing\")\n\n$d.Save()\n"}
